# Insert a new weekly record at row 14 ("Fruta / hortaliza, semanal"):
# every existing data row from 14 downward shifts down by one (14->15, ...,
# 117->118), and the freshly opened row 14 is populated with the new
# observation (week of 2021-12-07, $/saco 25 kilos, Región del Maule).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 14..117 down to 15..118, leaving row 14 blank (formatting carries
# over from the row above, same as a manual "Insert Row" in Excel).
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "Macroferia Regional de Talca"
$ws.Range("C14").Value = "Maule"
$ws.Range("D14").Value = 44537
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 100112031
$ws.Range("G14").Value = "Poroto verde"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 22000
$ws.Range("L14").Value = 22000
$ws.Range("M14").Value = 22000
$ws.Range("N14").Value = '$/saco 25 kilos'
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 880
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
